$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Range("E1").Value = "Uitleg"
$ws.Range("E2").Value = "Primary ID met sleutel"
$ws.Range("E3").Value = "Team naam"
$ws.Range("E4").Value = "Aantal spelers per team"
$ws.Range("E5").Value = "Aantal doelpunten per team"
$ws.Range("E6").Value = "Aantal keer gewonnen"
$ws.Range("E7").Value = "Aantal keer verloren"
$ws.Range("E9").Value = "Uitleg"
$ws.Range("E10").Value = "Wedstrijd ID"
$ws.Range("E11").Value = "Team ID, van tbl_teams"
$ws.Range("E12").Value = "Wedstrijd tijd"
$ws.Range("E13").Value = "De zaal van de wedstrijd"
$ws.Range("E14").Value = "Aantal doelpunten per wedstrijd"
$ws.Range("E16").Value = "Uitleg"
$ws.Range("E17").Value = "Speler ID's"
$ws.Range("E18").Value = "Speler namen"
$ws.Range("E19").Value = "Speler leeftijden"
$ws.Range("E20").Value = "Speler geslacht"
$ws.Range("E21").Value = "Teams ID van _tbl_teams"
$ws.Range("E22").Value = "Aantal doelpunten per speler"
$ws.Columns.Item(5).EntireColumn.AutoFit()
Write-Output ($ws.Columns.Item(5).ColumnWidth)
